$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D32").Value = "Logistic Regression은 왜 linear classification 인가?"
$ws.Range("E32").Value = "https://dodonam.tistory.com/392"

$ws.Range("D50").Value = "단백질 접힘"
$ws.Range("E50").Value = "http://incredible.egloos.com/7573055"

$ws.Range("D51").Value = "[금융데이터분석] 2022년 11월 3일 기준 삼성전자의 주가는 가치에 비해 싼가? 비싼가? (S-RIM 활용)"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EA%B8%88%EC%9C%B5%EB%8D%B0%EC%9D%B4%ED%84%B0%EB%B6%84%EC%84%9D-2022%EB%85%84-11%EC%9B%94-3%EC%9D%BC-%EA%B8%B0%EC%A4%80-%EC%82%BC%EC%84%B1%EC%A0%84%EC%9E%90%EC%9D%98-%EC%A3%BC%EA%B0%80%EB%8A%94-%EA%B0%80%EC%B9%98%EC%97%90-%EB%B9%84%ED%95%B4-%EC%8B%BC%EA%B0%80-%EB%B9%84%EC%8B%BC%EA%B0%80-S-RIM-%ED%99%9C%EC%9A%A9"
